$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 with the same style as H1
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J columns for each data row
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 29) {
        $ws.Cells.Item($r, 9).Value2 = 6
        $ws.Cells.Item($r, 10).Value2 = 7
    } else {
        $ws.Cells.Item($r, 9).Value2 = 1
        $ws.Cells.Item($r, 10).Value2 = $hVal
    }
}
